# msz - 3./4. smoke test + inheritance page and 2. dialog
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-fit columns for the new, longer content that is about to be added ---
$ws.Columns.Item(1).ColumnWidth = 53.05338541666667
$ws.Columns.Item(3).ColumnWidth = 43.05338541666667
$ws.Columns.Item(4).ColumnWidth = 31.498697916666668
$ws.Columns.Item(5).ColumnWidth = 43.05338541666667
$ws.Columns.Item(7).ColumnWidth = 43.05338541666667

# --- Append 4 new rows (15-18), mirroring the existing "<SET>"/"<NOP>" pattern ---
# Cell writes are ordered so brand-new shared-string entries are introduced in
# the same sequence the original workbook uses.

# Row 15: Vehicle Data page -> FillPage
$ws.Cells.Item(15, 1).Value = "104_MotorcycleInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Cells.Item(15, 2).Value = "<SET>"
$ws.Cells.Item(15, 3).Value = "104_MotorcycleInsurance_001_SmokeTest_FillPage"
$ws.Cells.Item(15, 8).Value = "<NOP>"

# Row 16: Insurant Data page -> DefaultInsurance FillPage (text-formatted)
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "DefaultInsurance_SmokeTest_FillPage"
$ws.Cells.Item(16, 1).Value = "104_MotorcycleInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Cells.Item(16, 2).Value = "<SET>"
$ws.Cells.Item(16, 8).Value = "<NOP>"

# Row 17: Product Data page -> FillPage (text-formatted, including col A)
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "104_MotorcycleInsurance_001_SmokeTest_FillPageProductData"
$ws.Cells.Item(17, 2).Value = "<SET>"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "104_MotorcycleInsurance_001_SmokeTest_FillPage"
$ws.Cells.Item(17, 8).Value = "<NOP>"

# Row 18: Send Quote page -> DefaultInsurance FillPage
$ws.Cells.Item(18, 1).Value = "104_MotorcycleInsurance_001_SmokeTest_FillPageSendQuote"
$ws.Cells.Item(18, 2).Value = "<SET>"
$ws.Cells.Item(18, 7).Value = "DefaultInsurance_SmokeTest_FillPage"
$ws.Cells.Item(18, 8).Value = "<NOP>"

# --- Selection moved to the last entered cell ---
$ws.Range("G18").Select() | Out-Null

# --- Move / resize the screenshot picture lower & narrower on the sheet ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 3.0
$shp.Top = 320.4
$shp.Width = 1124.85937007874
$shp.Height = 615.6730708661418
